$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row rename ---
$ws.Range('A1').Value = 'mx_state'
$ws.Range('B1').Value = 'mx_municipality'
$ws.Range('C1').Value = 'n_matriculas'
$ws.Range('D1').Value = 'pct_matriculas'

# --- A189 special case: strip embedded CR/LF + rename connector word ---
$ws.Range('A189').Value = 'Estado De México'

# --- Capitalize Spanish connector words (de/del/el/la/las/los/y) in state/municipality names ---
$nameFixes = @(
    @('B7', 'Pabellón De Arteaga'),
    @('B8', 'Rincón De Romos'),
    @('B9', 'San Francisco De Los Romo'),
    @('B10', 'San José De Gracia'),
    @('B14', 'Playas De Rosarito'),
    @('B32', 'Amatenango De La Frontera'),
    @('B33', 'Amatenango Del Valle'),
    @('B36', 'Bejucal De Ocampo'),
    @('B42', 'Chiapa De Corzo'),
    @('B45', 'Comitán De Domínguez'),
    @('B71', 'San Cristóbal De Las Casas'),
    @('B96', 'Guadalupe Y Calvo'),
    @('B98', 'Hidalgo Del Parral'),
    @('B108', 'San Francisco De Borja'),
    @('B109', 'San Francisco Del Oro'),
    @('B126', 'San Juan De Sabinas'),
    @('B139', 'Villa De Álvarez'),
    @('A141', 'Ciudad De México'),
    @('B145', 'Cuajimalpa De Morelos'),
    @('B160', 'Coneto De Comonfort'),
    @('B171', 'Nombre De Dios'),
    @('B179', 'San Juan Del Río'),
    @('B180', 'San Pedro Del Gallo'),
    @('B189', 'Acambay De Ruíz Castañeda'),
    @('B191', 'Almoloya De Alquisiras'),
    @('B192', 'Almoloya De Juárez'),
    @('B196', 'Atizapán De Zaragoza'),
    @('B203', 'Chapa De Mota'),
    @('B207', 'Coacalco De Berriozábal'),
    @('B212', 'Ecatepec De Morelos'),
    @('B217', 'Ixtapan De La Sal'),
    @('B218', 'Ixtapan Del Oro'),
    @('B229', 'Naucalpan De Juárez'),
    @('B237', 'San Felipe Del Progreso'),
    @('B248', 'Tenango Del Valle'),
    @('B258', 'Tlalnepantla De Baz'),
    @('B264', 'Valle De Chalco Solidaridad'),
    @('B265', 'Villa De Allende'),
    @('B266', 'Villa Del Carbón'),
    @('B276', 'San Miguel De Allende'),
    @('B277', 'Apaseo El Alto'),
    @('B285', 'Dolores Hidalgo Cuna De La Independencia Nacional'),
    @('B289', 'Jaral Del Progreso'),
    @('B295', 'Purísima Del Rincón'),
    @('B300', 'San Francisco Del Rincón'),
    @('B302', 'San Luis De La Paz'),
    @('B303', 'Santa Cruz De Juventino Rosas'),
    @('B304', 'Silao De La Victoria'),
    @('B309', 'Valle De Santiago'),
    @('B315', 'Acapulco De Juárez'),
    @('B317', 'Ajuchitlán Del Progreso'),
    @('B318', 'Alcozauca De Guerrero'),
    @('B322', 'Atenango Del Río'),
    @('B324', 'Atoyac De Álvarez'),
    @('B325', 'Ayutla De Los Libres'),
    @('B328', 'Buenavista De Cuéllar'),
    @('B329', 'Chilapa De Álvarez'),
    @('B330', 'Chilpancingo De Los Bravo'),
    @('B331', 'Coahuayutla De José María Izazaga'),
    @('B336', 'Coyuca De Benítez'),
    @('B337', 'Coyuca De Catalán'),
    @('B341', 'Cuetzala Del Progreso'),
    @('B342', 'Cutzamala De Pinzón'),
    @('B348', 'Huitzuco De Los Figueroa'),
    @('B349', 'Iguala De La Independencia'),
    @('B351', 'Ixcateopan De Cuauhtémoc'),
    @('B354', 'La Unión De Isidoro Montes De Oca'),
    @('B357', 'Mártir De Cuilapan'),
    @('B370', 'Taxco De Alarcón'),
    @('B372', 'Técpan De Galeana'),
    @('B374', 'Tepecoacuilco De Trujano'),
    @('B376', 'Tixtla De Guerrero'),
    @('B380', 'Tlalixtaquilla De Maldonado'),
    @('B381', 'Tlapa De Comonfort'),
    @('B395', 'Atotonilco El Grande'),
    @('B400', 'Cuautepec De Hinojosa'),
    @('B405', 'Huasca De Ocampo'),
    @('B408', 'Huejutla De Reyes'),
    @('B411', 'Jacala De Ledezma'),
    @('B417', 'Mineral Del Monte'),
    @('B418', 'Mixquiahuala De Juárez'),
    @('B419', 'Molango De Escamilla'),
    @('B420', 'Nopala De Villagrán'),
    @('B421', 'Omitlán De Juárez'),
    @('B422', 'Pachuca De Soto'),
    @('B423', 'Progreso De Obregón'),
    @('B428', 'Santiago De Anaya'),
    @('B429', 'Santiago Tulantepec De Lugo Guerrero'),
    @('B432', 'Tenango De Doria'),
    @('B434', 'Tepehuacán De Guerrero'),
    @('B435', 'Tepeji Del Río De Ocampo'),
    @('B438', 'Tezontepec De Aldama'),
    @('B441', 'Tula De Allende'),
    @('B442', 'Tulancingo De Bravo'),
    @('B443', 'Zacualtipán De Ángeles'),
    @('B448', 'Acatlán De Juárez'),
    @('B449', 'Ahualulco De Mercado'),
    @('B453', 'Atemajac De Brizuela'),
    @('B456', 'Atotonilco El Alto'),
    @('B458', 'Autlán De Navarro'),
    @('B464', 'Cañadas De Obregón'),
    @('B470', 'Concepción De Buenos Aires'),
    @('B471', 'Cuautitlán De García Barragán'),
    @('B479', 'Encarnación De Díaz'),
    @('B484', 'Huejuquilla El Alto'),
    @('B485', 'Ixtlahuacán De Los Membrillos'),
    @('B486', 'Ixtlahuacán Del Río'),
    @('B490', 'Jilotlán De Los Dolores'),
    @('B496', 'Lagos De Moreno'),
    @('B502', 'Ojuelos De Jalisco'),
    @('B507', 'San Cristóbal De La Barranca'),
    @('B509', 'San Juan De Los Lagos'),
    @('B510', 'San Juanito De Escobedo'),
    @('B512', 'San Martín De Bolaños'),
    @('B514', 'San Miguel El Alto'),
    @('B515', 'San Sebastián Del Oeste'),
    @('B518', 'Talpa De Allende'),
    @('B519', 'Tamazula De Gordiano'),
    @('B524', 'Teocuitatlán De Corona'),
    @('B525', 'Tepatitlán De Morelos'),
    @('B528', 'Tizapán El Alto'),
    @('B529', 'Tlajomulco De Zúñiga'),
    @('B541', 'Unión De San Antonio'),
    @('B542', 'Unión De Tula'),
    @('B543', 'Valle De Guadalupe'),
    @('B544', 'Valle De Juárez'),
    @('B548', 'Yahualica De González Gallo'),
    @('B549', 'Zacoalco De Torres'),
    @('B552', 'Zapotitlán De Vadillo'),
    @('B553', 'Zapotlán Del Rey'),
    @('B554', 'Zapotlán El Grande'),
    @('B578', 'Coalcomán De Vázquez Pallares'),
    @('B580', 'Cojumatlán De Régules'),
    @('B642', 'Tiquicheo De Nicolás Romero'),
    @('B667', 'Coatlán Del Río'),
    @('B677', 'Puente De Ixtla'),
    @('B681', 'Tetela Del Volcán'),
    @('B682', 'Tlaltizapán De Zapata'),
    @('B692', 'Bahía De Banderas'),
    @('B696', 'Ixtlán Del Río'),
    @('B703', 'Santa María Del Oro'),
    @('B718', 'San Nicolás De Los Garza'),
    @('B721', 'Acatlán De Pérez Figueroa'),
    @('B725', 'Ayoquezco De Aldama'),
    @('B729', 'Chalcatongo De Hidalgo'),
    @('B731', 'Coicoyán De Las Flores'),
    @('B734', 'Constancia Del Rosario'),
    @('B736', 'Cuilápam De Guerrero'),
    @('B737', 'Guadalupe De Ramírez'),
    @('B738', 'Guevea De Humboldt'),
    @('B739', 'Heroica Ciudad De Ejutla De Crespo'),
    @('B740', 'Heroica Ciudad De Huajuapan De León'),
    @('B741', 'Heroica Ciudad De Tlaxiaco'),
    @('B742', 'Huautla De Jiménez'),
    @('B744', 'Ixtlán De Juárez'),
    @('B745', 'Heroica Ciudad De Juchitán De Zaragoza'),
    @('B754', 'Magdalena Yodocono De Porfirio Díaz'),
    @('B755', 'Mariscala De Juárez'),
    @('B758', 'Miahuatlán De Porfirio Díaz'),
    @('B759', 'Mixistlán De La Reforma'),
    @('B761', 'Nejapa De Madero'),
    @('B763', 'Oaxaca De Juárez'),
    @('B764', 'Ocotlán De Morelos'),
    @('B765', 'Putla Villa De Guerrero'),
    @('B766', 'Reforma De Pineda'),
    @('B776', 'San Antonio De La Cal'),
    @('B782', 'San Baltazar Yatzachi El Bajo'),
    @('B790', 'San Felipe Jalapa De Díaz'),
    @('B809', 'San José Del Peñasco'),
    @('B810', 'San José Del Progreso'),
    @('B819', 'San Juan Del Estado'),
    @('B820', 'San Juan Del Río'),
    @('B852', 'San Mateo Del Mar'),
    @('B858', 'San Miguel El Grande'),
    @('B878', 'San Pedro El Alto'),
    @('B890', 'San Pedro Y San Pablo Teposcolula'),
    @('B895', 'Santa Ana Del Valle'),
    @('B903', 'Santa Cruz De Bravo'),
    @('B908', 'Santa Cruz Tacache De Mina'),
    @('B911', 'Santa Inés Del Monte'),
    @('B924', 'Santa María Jalapa Del Marqués'),
    @('B965', 'Santo Domingo De Morelos'),
    @('B978', 'Sitio De Xitlapehua'),
    @('B980', 'Tamazulápam Del Espíritu Santo'),
    @('B981', 'Tanetze De Zaragoza'),
    @('B982', 'Tataltepec De Valdés'),
    @('B983', 'Teococuilco De Marcos Pérez'),
    @('B984', 'Teotitlán De Flores Magón'),
    @('B985', 'Teotitlán Del Valle'),
    @('B987', 'Tlacolula De Matamoros'),
    @('B988', 'Tlalixtac De Cabrera'),
    @('B989', 'Totontepec Villa De Morelos'),
    @('B992', 'Villa De Chilapa De Díaz'),
    @('B993', 'Villa De Etla'),
    @('B994', 'Villa De Tamazulápam Del Progreso'),
    @('B995', 'Villa De Tututepec De Melchor Ocampo'),
    @('B996', 'Villa De Zaachila'),
    @('B998', 'Villa Sola De Vega'),
    @('B1003', 'Zimatlán De Álvarez'),
    @('B1021', 'Chalchicomula De Sesma'),
    @('B1035', 'Cuayuca De Andrade'),
    @('B1036', 'Cuetzalan Del Progreso'),
    @('B1048', 'Huehuetlán El Grande'),
    @('B1051', 'Ixcamilpa De Guerrero'),
    @('B1054', 'Izúcar De Matamoros'),
    @('B1063', 'Los Reyes De Juárez'),
    @('B1071', 'Palmar De Bravo'),
    @('B1080', 'San Diego La Mesa Tochimiltzingo'),
    @('B1092', 'San Salvador El Seco'),
    @('B1093', 'San Salvador El Verde'),
    @('B1098', 'Tecali De Herrera'),
    @('B1105', 'Tepanco De López'),
    @('B1106', 'Tepatlaxco De Hidalgo'),
    @('B1110', 'Tepexi De Rodríguez'),
    @('B1112', 'Tetela De Ocampo'),
    @('B1117', 'Tlacotepec De Benito Juárez'),
    @('B1141', 'Amealco De Bonfil'),
    @('B1142', 'Cadereyta De Montes'),
    @('B1147', 'Jalpan De Serra'),
    @('B1149', 'Pinal De Amoles'),
    @('B1152', 'San Juan Del Río'),
    @('B1163', 'Armadillo De Los Infante'),
    @('B1168', 'Ciudad Del Maíz'),
    @('B1173', 'Mexquitic De Carmona'),
    @('B1180', 'Soledad De Graciano Sánchez'),
    @('B1186', 'Villa De Arista'),
    @('B1187', 'Villa De Arriaga'),
    @('B1188', 'Villa De Ramos'),
    @('B1189', 'Villa De Reyes'),
    @('B1224', 'Nacozari De García'),
    @('B1262', 'Apetatitlán De Antonio Carvajal'),
    @('B1270', 'Ixtacuixtla De Mariano Matamoros'),
    @('B1290', 'Amatlán De Los Reyes'),
    @('B1310', 'Cosamaloapan De Carpio'),
    @('B1325', 'Hueyapan De Ocampo'),
    @('B1326', 'Ignacio De La Llave'),
    @('B1329', 'Ixhuatlán De Madero'),
    @('B1330', 'Ixhuatlán Del Café'),
    @('B1337', 'Juchique De Ferrer'),
    @('B1341', 'Lerdo De Tejada'),
    @('B1345', 'Martínez De La Torre'),
    @('B1348', 'Mixtla De Altamirano'),
    @('B1357', 'Paso Del Macho'),
    @('B1359', 'Poza Rica De Hidalgo'),
    @('B1365', 'Sayula De Alemán'),
    @('B1368', 'Soledad De Doblado'),
    @('B1391', 'Vega De Alatorre'),
    @('B1399', 'Zontecomatlán De López Y Fuentes'),
    @('B1434', 'El Plateado De Joaquín Amaro'),
    @('B1442', 'Jiménez Del Teul'),
    @('B1446', 'Mezquital Del Oro'),
    @('B1451', 'Moyahua De Estrada'),
    @('B1452', 'Nochistlán De Mejía'),
    @('B1453', 'Noria De Ángeles'),
    @('B1461', 'Teúl De González Ortega'),
    @('B1462', 'Tlaltenango De Sánchez Román'),
    @('B1463', 'Trinidad García De La Cadena'),
    @('B1465', 'Villa De Cos'),
)
foreach ($fix in $nameFixes) {
    $ws.Range($fix[0]).Value = $fix[1]
}

# --- Floating point cache refresh on percentage column (D) ---
$pctFixes = @(
    @('D4', 0.0009997857601942439),
    @('D17', 0.009497964721845321),
    @('D21', 0.0009997857601942439),
    @('D156', 0.0009997857601942439),
    @('D298', 0.0009997857601942439),
    @('D323', 0.0009997857601942439),
    @('D333', 0.0009997857601942439),
    @('D344', 0.0009997857601942439),
    @('D345', 0.0009997857601942439),
    @('D410', 0.0009997857601942439),
    @('D424', 0.0009997857601942439),
    @('D496', 0.0009997857601942439),
    @('D545', 0.0009997857601942439),
    @('D586', 0.0009997857601942439),
    @('D638', 0.0009997857601942439),
    @('D660', 0.0009997857601942439),
    @('D740', 0.009283724916089408),
    @('D755', 0.0009997857601942439),
    @('D791', 0.0009997857601942439),
    @('D976', 0.0009997857601942439),
    @('D1197', 0.0009997857601942439),
    @('D1412', 0.0009997857601942439),
)
foreach ($fix in $pctFixes) {
    $ws.Range($fix[0]).Value = $fix[1]
}

# --- Remove trailing footnote rows (now out of the data range) ---
$ws.Rows('1474:1478').Delete()

